$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Keep a pristine, unformatted "normal" cell to borrow clean formatting from
# (it is never itself modified by this script), so that cells whose literal
# text happens to be "true"/"false" don't end up stuck with the quote-prefix
# style that typing a leading apostrophe leaves behind.
$fmtSource = $ws.Range("B9")

# Experimental = true
$ws.Range("B7").Value = "'true"
$fmtSource.Copy()
$ws.Range("B7").PasteSpecial(-4122)

# Date updated
$ws.Range("B8").Value = "2025-07-21T12:46:15+00:00"

# Case Sensitive stays true
$ws.Range("B15").Value = "'true"
$fmtSource.Copy()
$ws.Range("B15").PasteSpecial(-4122)

# Compositional = false
$ws.Range("B18").Value = "'false"
$fmtSource.Copy()
$ws.Range("B18").PasteSpecial(-4122)

$excel.CutCopyMode = $false
